# Add a "Save" column (column H) to the s_vals sheet, mirroring the
# header style used by the existing columns (B1:G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (G1, style index 1:
# bold, centered, bordered) onto the new header cell H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data column values.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
